$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.735474467277527
$ws.Range("B1").Value = 2.497875928878784
$ws.Range("C1").Value = 2.579309701919556
$ws.Range("D1").Value = 2.953764915466309
$ws.Range("E1").Value = 3.666409254074097
